$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G9 changes from shared string "0.01" to shared string "0.0"
$ws.Range("G9").Value = "0.0"

# H9 changes from numeric 0 to shared string "0.1"
$ws.Range("H9").Value = "0.1"

# Update the active selection to H9, matching the saved view state
$ws.Range("H9").Select()
